# Add some temporary interests for Bastien (Outdoors, Culinary Arts, Biotechnology)
# in the "Interests" section of the Author form (rows 35-37 of column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "Outdoors"
$ws.Range("A36").Value = "Culinary Arts"
$ws.Range("A37").Value = "Biotechnology"

# Leave the cursor where the author finished typing, matching the scrolled
# view captured in the saved workbook (topLeftCell A27, active cell B38).
[void]$ws.Range("B38").Select()
